# Update the 农村居民消费 (rural household consumption) table:
#  - drop the oldest four year-rows (2000, 2002, 2005, 2007)
#  - keep 2010 / 2012 / 2015 / 2017 (they shift up into rows 2-5)
#  - append a new 2020 row with fresh figures

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the four obsolete year rows (2000年, 2002年, 2005年, 2007年).
# This shifts the remaining data rows (2010/2012/2015/2017, previously
# rows 6-9) up to rows 2-5, exactly matching the target layout.
$ws.Range("A2:A5").EntireRow.Delete()

# Give the new row's label cell (A6) the same style as the other year
# labels (bold/bordered/centered) before writing into it.
$ws.Range("A5").Copy()
$ws.Range("A6").PasteSpecial(-4122)

# New 2020年 row of figures.
$ws.Range("A6").Value = "2020年"
$ws.Range("B6").Value = 1188015.95895082
$ws.Range("C6").Value = 155878057.103352
$ws.Range("D6").Value = 109745451.401855
$ws.Range("E6").Value = 816642871.617192
$ws.Range("J6").Value = 54242995.1094221
$ws.Range("K6").Value = 26404008.7322014
$ws.Range("L6").Value = 16110906.1592875
$ws.Range("M6").Value = 23013748.4379359
$ws.Range("O6").Value = 753293.483239878
$ws.Range("P6").Value = 665462.948125652
$ws.Range("R6").Value = 1560823.70183441
$ws.Range("S6").Value = 156188434.737491
# F6, G6, H6, I6, N6, Q6 have no reported value for 2020年 and are left blank.
